# Regenerate the handoff report: the e2e test file was renamed (new GUID)
# and new Xliff hand-off/hand-back files were produced, so every cell that
# referenced the old identifiers needs to point at the new ones, and the
# "Ready for handoff" timestamps move forward a few seconds.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldGuid = "1f25c435-9197-4cbe-b0d3-b3b4c21d5293"
$newGuid = "35fcf230-f3f3-499f-8195-5edfd46dc5d4"

$oldHash = "2367c8074114a1edf9a9a0559dbd2f456a7c6ff3"
$newHash = "54bad4e69842081a8f0290b0392a202a75f3a2e7"

# The hyperlinks' underlying target address (commit blob URL) is left as-is
# by this change -- only the visible text changes to the new file name.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/35dcfe3a9bccf56cad9b9bc2aafd37ddc1337006/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview.Range("A2").Value() = ($newGuid + ".md")
$overview.Range("G2").Value() = "2016-08-28 22:57:20"

$bLink = $overview.Range("B2")
$bLink.Hyperlinks.Delete()
$overview.Hyperlinks.Add($bLink, $linkAddress, [System.Type]::Missing, [System.Type]::Missing, ("e2e\" + $newGuid + ".md")) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn.Range("G2").Value() = ($newGuid + "." + $newHash + ".zh-cn.xlf")
$zhcn.Range("H2").Value() = "2016-08-28 22:57:15"

$zhLink = $zhcn.Range("A2")
$zhLink.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhLink, $linkAddress, [System.Type]::Missing, [System.Type]::Missing, ($newGuid + ".md")) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede.Range("G2").Value() = ($newGuid + "." + $newHash + ".de-de.xlf")
$dede.Range("H2").Value() = "2016-08-28 22:57:20"

$deLink = $dede.Range("A2")
$deLink.Hyperlinks.Delete()
$dede.Hyperlinks.Add($deLink, $linkAddress, [System.Type]::Missing, [System.Type]::Missing, ($newGuid + ".md")) | Out-Null

# ---------------------------------------------------------------------
# Column A on every sheet is very slightly narrower in the regenerated
# report (best-effort: COM ColumnWidth only supports 1/6-character
# granularity, so we pick the closest reachable width to 39.6252049037388).
# ---------------------------------------------------------------------
$overview.Columns.Item(1).ColumnWidth = 38.8333333333
$zhcn.Columns.Item(1).ColumnWidth = 38.8333333333
$dede.Columns.Item(1).ColumnWidth = 38.8333333333
